# Changement de theme fonctionnel
# L'utilisateur peut changer le theme de l'application.
#
# This script mirrors the author's edits:
#   - New journal entries added to "Iteration #3" (rows 14-15-16) and
#     "Iteration #2" (row 24), each with a date / description / hours.
#   - The active sheet moves from "Iteration #2" to "Iteration #3".

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Iteration #2")
$ws3 = $wb.Worksheets.Item("Iteration #3")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# "Iteration #3" rows 14-15 : new work-log entries
# (written first so the new shared strings land at indices 59 and 60,
#  matching the order they were originally authored in)
# ---------------------------------------------------------------------

# Row 14 : 2017-03-27, "ajout d'une nouvelle table ...", 3h
$ws2.Cells.Item(15, 1).Copy()
$ws3.Cells.Item(14, 1).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(14, 1).Value = 42821

$ws2.Cells.Item(14, 2).Copy()
$ws3.Cells.Item(14, 2).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(14, 2).Value = "ajout d'une nouvelle table dans ma BD SqLite pour stocker mon theme, reglage de bug"

$ws3.Cells.Item(14, 3).Value = 3

$ws3.Rows.Item(14).RowHeight = 28.8

# Row 15 : 2017-03-28, "réglage du bug lors de la creation ...", 1.5h
$ws2.Cells.Item(16, 1).Copy()
$ws3.Cells.Item(15, 1).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(15, 1).Value = 42822

$ws2.Cells.Item(19, 2).Copy()
$ws3.Cells.Item(15, 2).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(15, 2).Value = "réglage du bug lors de la creation de mes tables dans SQLLite"

$ws3.Cells.Item(15, 3).Value = 1.5

# ---------------------------------------------------------------------
# "Iteration #2" row 24 : new work-log entry (shared string index 61)
# ---------------------------------------------------------------------

$ws2.Cells.Item(15, 1).Copy()
$ws2.Cells.Item(24, 1).PasteSpecial($xlPasteFormats)
$ws2.Cells.Item(24, 1).Value = 42814

$ws2.Cells.Item(24, 2).Value = "ajout d'un scrool view dans l'activite des themes (journée des présentations)"

$ws2.Cells.Item(24, 3).Value = 3

# the now-unused placeholder cells below it lose their formatting
$ws2.Cells.Item(25, 1).Clear()
$ws2.Cells.Item(26, 1).Clear()

# ---------------------------------------------------------------------
# "Iteration #3" row 16 : new work-log entry (shared string index 62)
# ---------------------------------------------------------------------

$ws2.Cells.Item(17, 1).Copy()
$ws3.Cells.Item(16, 1).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(16, 1).Value = 42828

$ws2.Cells.Item(14, 2).Copy()
$ws3.Cells.Item(16, 2).PasteSpecial($xlPasteFormats)
$ws3.Cells.Item(16, 2).Value = " Tests sur le fonctionnement de la sauvegarde du theme actif"

$ws3.Cells.Item(16, 3).Value = 1.5

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------

$ws2.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
$ws2.Range("B28").Select() | Out-Null

$ws3.Activate() | Out-Null
$win3 = $excel.ActiveWindow
$win3.ScrollRow = 7
$win3.ScrollColumn = 1
$ws3.Range("C15").Select() | Out-Null
